$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log_Muestras")

$timestamps = @{
  2 = "2025-11-13T06:52:51.265673"
  3 = "2025-11-13T06:52:51.270441"
  4 = "2025-11-13T06:52:51.270441"
  5 = "2025-11-13T06:52:51.271444"
  6 = "2025-11-13T06:52:51.271752"
  7 = "2025-11-13T06:52:51.271752"
  8 = "2025-11-13T06:52:51.272267"
  9 = "2025-11-13T06:52:51.272267"
  10 = "2025-11-13T06:52:51.272267"
  11 = "2025-11-13T06:52:51.272267"
  12 = "2025-11-13T06:52:51.273283"
  13 = "2025-11-13T06:52:51.273283"
  14 = "2025-11-13T06:52:51.273283"
  15 = "2025-11-13T06:52:51.273283"
  16 = "2025-11-13T06:52:51.273283"
  17 = "2025-11-13T06:52:51.274286"
  18 = "2025-11-13T06:52:51.274286"
  19 = "2025-11-13T06:52:51.274286"
  20 = "2025-11-13T06:52:51.274286"
  21 = "2025-11-13T06:52:51.275285"
  22 = "2025-11-13T06:52:51.275285"
  23 = "2025-11-13T06:52:51.275285"
  24 = "2025-11-13T06:52:51.275285"
  25 = "2025-11-13T06:52:51.275285"
  26 = "2025-11-13T06:52:51.276287"
  27 = "2025-11-13T06:52:51.276287"
  28 = "2025-11-13T06:52:51.276287"
  29 = "2025-11-13T06:52:51.277284"
  30 = "2025-11-13T06:52:51.277284"
  31 = "2025-11-13T06:52:51.277284"
  32 = "2025-11-13T06:52:51.277284"
  33 = "2025-11-13T06:52:51.278285"
  34 = "2025-11-13T06:52:51.278285"
  35 = "2025-11-13T06:52:51.278285"
  36 = "2025-11-13T06:52:51.278285"
  37 = "2025-11-13T06:52:51.279285"
  38 = "2025-11-13T06:52:51.279285"
  39 = "2025-11-13T06:52:51.279285"
  40 = "2025-11-13T06:52:51.279285"
  41 = "2025-11-13T06:52:51.279285"
  42 = "2025-11-13T06:52:51.280286"
  43 = "2025-11-13T06:52:51.281284"
  44 = "2025-11-13T06:52:51.282283"
  45 = "2025-11-13T06:52:51.282283"
  46 = "2025-11-13T06:52:51.621316"
  47 = "2025-11-13T06:52:51.621316"
  48 = "2025-11-13T06:52:51.622312"
  49 = "2025-11-13T06:52:51.622312"
  50 = "2025-11-13T06:52:51.622312"
  51 = "2025-11-13T06:52:51.622312"
  52 = "2025-11-13T06:52:51.622312"
  53 = "2025-11-13T06:52:51.622312"
  54 = "2025-11-13T06:52:51.622312"
  55 = "2025-11-13T06:52:51.622312"
  56 = "2025-11-13T06:52:51.622312"
  57 = "2025-11-13T06:52:51.623311"
  58 = "2025-11-13T06:52:51.623311"
  59 = "2025-11-13T06:52:51.623311"
  60 = "2025-11-13T06:52:51.623311"
  61 = "2025-11-13T06:52:51.623311"
  62 = "2025-11-13T06:52:51.623311"
  63 = "2025-11-13T06:52:51.623311"
  64 = "2025-11-13T06:52:51.623311"
  65 = "2025-11-13T06:52:51.623311"
  66 = "2025-11-13T06:52:51.623311"
  67 = "2025-11-13T06:52:51.623311"
  68 = "2025-11-13T06:52:51.624311"
  69 = "2025-11-13T06:52:51.624311"
  70 = "2025-11-13T06:52:51.624311"
  71 = "2025-11-13T06:52:51.624311"
  72 = "2025-11-13T06:52:51.624311"
  73 = "2025-11-13T06:52:51.625312"
  74 = "2025-11-13T06:52:51.625312"
  75 = "2025-11-13T06:52:51.820158"
  76 = "2025-11-13T06:52:51.820158"
  77 = "2025-11-13T06:52:51.820158"
  78 = "2025-11-13T06:52:51.820158"
  79 = "2025-11-13T06:52:51.820158"
  80 = "2025-11-13T06:52:51.820158"
  81 = "2025-11-13T06:52:51.820158"
  82 = "2025-11-13T06:52:51.821166"
  83 = "2025-11-13T06:52:51.821166"
  84 = "2025-11-13T06:52:51.821166"
  85 = "2025-11-13T06:52:51.821166"
  86 = "2025-11-13T06:52:51.821166"
  87 = "2025-11-13T06:52:51.821166"
  88 = "2025-11-13T06:52:51.821166"
  89 = "2025-11-13T06:52:51.821166"
  90 = "2025-11-13T06:52:51.821166"
  91 = "2025-11-13T06:52:51.821166"
  92 = "2025-11-13T06:52:51.822162"
  93 = "2025-11-13T06:52:51.822162"
  94 = "2025-11-13T06:52:51.822162"
  95 = "2025-11-13T06:52:51.822162"
  96 = "2025-11-13T06:52:51.822162"
  97 = "2025-11-13T06:52:51.822162"
  98 = "2025-11-13T06:52:51.822162"
  99 = "2025-11-13T06:52:51.822162"
  100 = "2025-11-13T06:52:51.822162"
  101 = "2025-11-13T06:52:51.822162"
  102 = "2025-11-13T06:52:51.822162"
}

foreach ($row in $timestamps.Keys) {
  $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}
